# Update "想去人数" (F column) and "最低票价" (G column) figures for the
# 展览 (Exhibition) and 全部类型 (All types) worksheets, which share the
# same underlying data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F7").Value = 109
    $ws.Range("F11").Value = 36
    $ws.Range("F14").Value = 334
    $ws.Range("F17").Value = 381
    $ws.Range("G19").Value = 35
    $ws.Range("F21").Value = 37
    $ws.Range("F22").Value = 911
    $ws.Range("F23").Value = 2697
    $ws.Range("F31").Value = 380
    $ws.Range("F33").Value = 591
    $ws.Range("F34").Value = 422
}
